$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-13 (Generation 0-11): set Fitness (column C) to 7865
$ws.Range("C2:C13").Value = 7865

# Rows 14-252 (Generation 12-250): set Fitness (column C) to 7293
$ws.Range("C14:C252").Value = 7293
